# The deck's single slide master currently uses the "Integral" theme
# (ppt/theme/theme2.xml) while the notes master uses a plain
# "Office Theme" (ppt/theme/theme1.xml). The authored change swaps the
# two themes' content: the master's theme becomes the "Office Theme"
# colour scheme, and vice versa. Re-colour the master's theme (the one
# that drives every slide) to the Office Theme palette via the Theme
# Color Scheme, which is the supported COM surface for editing a
# presentation's live theme colours.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# msoThemeColorDark1=1, Light1=2, Dark2=3, Light2=4,
# Accent1..6=5..10, Hyperlink=11, FollowedHyperlink=12
# Target values: the "Office Theme" colour scheme (the palette that
# previously lived in ppt/theme/theme1.xml).
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
